$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version: 5.0.0 -> 6.0.0
$ws.Range("B3").Value = "6.0.0"

# Date: 2021-12-16T17:36:56+00:00 -> 2022-01-21T20:46:54+00:00
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value (was blank) -> Alvearie Team
$ws.Range("B9").Value = "Alvearie Team"

# Remove the duplicated "Contact" / "No display for ContactDetail" row (old row 11)
$ws.Rows.Item(11).Delete()

# Old row 10 ("Contact" / "No display for ContactDetail") becomes "Jurisdiction" / "United States of America"
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# "Case Sensitive" row gets the text value "true" (kept as literal text, not boolean)
$ws.Range("F1").Value = "'true"
$ws.Range("F1").Copy()
$ws.Range("B14").PasteSpecial(-4163)
$ws.Range("F1").Clear()
$excel.CutCopyMode = 0
